$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "commit changes" (was B29, style s="4") up into B26
$ws.Range("B29").Cut($ws.Range("B26"))

# Move "End" (was B30, style s="5") up into B27
$ws.Range("B30").Cut($ws.Range("B27"))

# Fully clear the now-vacated cells B28:B30 (contents + formatting)
$ws.Range("B28:B30").Clear()

# Update the active selection to B30, matching the saved view state
$ws.Range("B30").Select()
